$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gearbox Tests")

$ws.Range("C4").Value = 731851.44965614588
$ws.Range("C5").Value = 128850.48535948661
$ws.Range("C6").Value = 3802.9230326870379
$ws.Range("C7").Value = 54455.593060061852
$ws.Range("C8").Value = 90
